$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Label" header in H1, matching the style of the other header cells (B1:G1) ---
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- Fill H2:H21 with the 0/1 class-label column ---
$labels = @{
    2  = 0;  3  = 0;  4  = 0;  5  = 0;  6  = 0
    7  = 1;  8  = 1;  9  = 1;  10 = 1;  11 = 1
    12 = 0;  13 = 0;  14 = 0;  15 = 0;  16 = 0
    17 = 1;  18 = 1;  19 = 1;  20 = 1;  21 = 1
}
foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}

# --- Refit values: small numeric updates to D/E/F from the refitted NCDEs ---
$ws.Range("D3").Value  = 0.4890206267493784
$ws.Range("E3").Value  = 0.4890206267493784

$ws.Range("D4").Value  = 0.3882675107359504
$ws.Range("E4").Value  = 0.3882675107359504

$ws.Range("D7").Value  = 0.3709167571845889
$ws.Range("E7").Value  = 0.6290832428154112

$ws.Range("D11").Value = 0.4197213079337474
$ws.Range("E11").Value = 0.5802786920662526
$ws.Range("F11").Value = 0.5812770128250122
